$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("SoCDTtiNTY")

# Update the LDVs passenger share (SoCDTtiNTY!B2): lifetime assumption
# dropped, raising the new-this-year share from 1/17 to 1/12.
$wsData.Activate()
$wsData.Range("B2").Value = 0.083333333333333329

# Move the selection on the SoCDTtiNTY sheet (it is no longer the active tab).
$wsData.Range("E41").Select()

# Make "About" the active/selected sheet (previously SoCDTtiNTY was active).
$wsAbout.Activate()
